$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Test Case"

# Select/format the Status column result cells (I11:I14): center them and
# box them with a thin border, matching the other test-case-manual styling.
$range = $ws.Range("I11:I14")
$range.HorizontalAlignment = -4108   # xlCenter
$range.VerticalAlignment = -4108     # xlCenter
$range.Borders.LineStyle = 1         # xlContinuous
$range.Borders.Weight = 2            # xlThin

# Leave the view scrolled to where the edited cells are, with them selected
$ws.Application.ActiveWindow.ScrollRow = 7
$sel = $ws.Range("I11:I14")
$sel.Select()
